$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 20
$ws.Range("B6").Value = 20

$ws.Range("D6").Select()
